# Auto-generated Excel COM-interop script applying the Halicarnassus_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1098
$ws.Range("I6").Value = 425.85715
$ws.Range("J6").Value = 2666.3333
$ws.Range("K6").Value = 1277.57145
$ws.Range("L6").Value = 7998.999899999999
$ws.Range("M6").Value = -1165.57145
$ws.Range("N6").Value = -8222.999899999999

# Row 58
$ws.Range("H58").Value = 1689.5714
$ws.Range("I58").Value = 109
$ws.Range("J58").Value = 2875
$ws.Range("K58").Value = 327
$ws.Range("L58").Value = 8625
$ws.Range("M58").Value = -177
$ws.Range("N58").Value = -8925

# Row 96
$ws.Range("H96").Value = 238.72728
$ws.Range("I96").Value = 238.72728
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 716.18184
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 656.81816
$ws.Range("N96").ClearContents()

# Row 100
$ws.Range("H100").Value = 5098
$ws.Range("I100").Value = 4001
$ws.Range("K100").Value = 4001
$ws.Range("M100").Value = -3460

# Row 125
$ws.Range("H125").Value = 2029.7142
$ws.Range("I125").Value = 1782.0667
$ws.Range("K125").Value = 16038.6003
$ws.Range("M125").Value = -13578.6003

# Row 135
$ws.Range("H135").Value = 1641.3846
$ws.Range("I135").Value = 1272.7778
$ws.Range("J135").Value = 2470.75
$ws.Range("K135").Value = 11455.0002
$ws.Range("L135").Value = 22236.75
$ws.Range("M135").Value = -8920.0002
$ws.Range("N135").Value = -27306.75

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 6135.8823
$ws.Range("I2").Value = 3675.6875
$ws.Range("J2").Value = 45499
$ws.Range("K2").Value = 3675.6875
$ws.Range("L2").Value = 45499
$ws.Range("M2").Value = -3562.6875
$ws.Range("N2").Value = -45725

# Row 61
$ws.Range("H61").Value = 6671.3335
$ws.Range("J61").Value = 10014
$ws.Range("L61").Value = 10014
$ws.Range("N61").Value = -10438

# Row 88
$ws.Range("H88").Value = 1629.625
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 1748.1428
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 1748.1428
$ws.Range("M88").Value = -394
$ws.Range("N88").Value = -2560.1428

# Row 91
$ws.Range("H91").Value = 1629.625
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 1748.1428
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 1748.1428
$ws.Range("M91").Value = 604
$ws.Range("N91").Value = -4556.1428

# Row 102
$ws.Range("H102").Value = 5383
$ws.Range("I102").Value = 2305
$ws.Range("K102").Value = 2305
$ws.Range("M102").Value = -683

# Row 116
$ws.Range("H116").Value = 6135.8823
$ws.Range("I116").Value = 3675.6875
$ws.Range("J116").Value = 45499
$ws.Range("K116").Value = 3675.6875
$ws.Range("L116").Value = 45499
$ws.Range("M116").Value = -1381.6875
$ws.Range("N116").Value = -50087

# Row 124
$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820

# Row 136
$ws.Range("H136").Value = 6671.3335
$ws.Range("J136").Value = 10014
$ws.Range("L136").Value = 30042
$ws.Range("N136").Value = -35142

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 6135.8823
$ws.Range("I3").Value = 3675.6875
$ws.Range("J3").Value = 45499
$ws.Range("K3").Value = 3675.6875
$ws.Range("L3").Value = 45499
$ws.Range("M3").Value = -3561.6875
$ws.Range("N3").Value = -45727

# Row 107
$ws.Range("H107").Value = 3427.4333
$ws.Range("I107").Value = 1091.2
$ws.Range("K107").Value = 1091.2
$ws.Range("M107").Value = 828.8

$ws = $wb.Worksheets.Item("CRP")
# Row 92
$ws.Range("H92").Value = 33749.5
$ws.Range("J92").Value = 38332.668
$ws.Range("L92").Value = 38332.668
$ws.Range("N92").Value = -43324.668

# Row 103
$ws.Range("H103").Value = 15262.5
$ws.Range("I103").Value = 15262.5
$ws.Range("K103").Value = 15262.5
$ws.Range("M103").Value = -14090.5

# Row 122
$ws.Range("H122").Value = 1303.4546
$ws.Range("I122").Value = 1174.2858
$ws.Range("J122").Value = 1529.5
$ws.Range("K122").Value = 3522.8574
$ws.Range("L122").Value = 4588.5
$ws.Range("M122").Value = -1072.8574
$ws.Range("N122").Value = -9488.5

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 358.18182
$ws.Range("I14").Value = 358.18182
$ws.Range("K14").Value = 1074.54546
$ws.Range("M14").Value = -901.54546

# Row 132
$ws.Range("H132").Value = 1942.7142
$ws.Range("J132").Value = 1942.7142
$ws.Range("L132").Value = 17484.4278
$ws.Range("N132").Value = -22544.4278

$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 20000642
$ws.Range("J7").Value = 24001540
$ws.Range("L7").Value = 24001540
$ws.Range("N7").Value = -24001764

# Row 8
$ws.Range("H8").Value = 20000642
$ws.Range("J8").Value = 24001540
$ws.Range("L8").Value = 24001540
$ws.Range("N8").Value = -24001818

# Row 102
$ws.Range("H102").Value = 1097.5
$ws.Range("I102").Value = 1097.5
$ws.Range("K102").Value = 1097.5
$ws.Range("M102").Value = 524.5

# Row 113
$ws.Range("H113").Value = 6529.7
$ws.Range("I113").Value = 4145.1816
$ws.Range("J113").Value = 9444.111000000001
$ws.Range("K113").Value = 4145.1816
$ws.Range("L113").Value = 9444.111000000001
$ws.Range("M113").Value = -1975.1816
$ws.Range("N113").Value = -13784.111

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 817.4286
$ws.Range("I16").Value = 753.8333
$ws.Range("J16").Value = 1199
$ws.Range("K16").Value = 753.8333
$ws.Range("L16").Value = 1199
$ws.Range("M16").Value = -583.8333
$ws.Range("N16").Value = -1539

# Row 22
$ws.Range("H22").Value = 1040.4
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590

# Row 27
$ws.Range("H27").Value = 1040.4
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214

# Row 30
$ws.Range("H30").Value = 459.6
$ws.Range("I30").Value = 459.6
$ws.Range("K30").Value = 459.6
$ws.Range("M30").Value = -351.6

# Row 40
$ws.Range("H40").Value = 5688.231
$ws.Range("I40").Value = 3192.7144
$ws.Range("J40").Value = 8599.666999999999
$ws.Range("K40").Value = 3192.7144
$ws.Range("L40").Value = 8599.666999999999
$ws.Range("M40").Value = -3056.7144
$ws.Range("N40").Value = -8871.666999999999

# Row 46
$ws.Range("H46").Value = 5204.8
$ws.Range("I46").Value = 3729.0908
$ws.Range("J46").Value = 6364.2856
$ws.Range("K46").Value = 3729.0908
$ws.Range("L46").Value = 6364.2856
$ws.Range("M46").Value = -3541.0908
$ws.Range("N46").Value = -6740.2856

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 502.06897
$ws.Range("I14").Value = 591.7646999999999
$ws.Range("J14").Value = 375
$ws.Range("K14").Value = 591.7646999999999
$ws.Range("L14").Value = 375
$ws.Range("M14").Value = -423.7646999999999
$ws.Range("N14").Value = -711

# Row 64
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496

# Row 67
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716

# Row 107
$ws.Range("H107").Value = 721
$ws.Range("I107").Value = 721
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2163
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -243
$ws.Range("N107").ClearContents()

# Row 122
$ws.Range("H122").Value = 1835
$ws.Range("I122").Value = 1835
$ws.Range("K122").Value = 5505
$ws.Range("M122").Value = -3055
